$d = $word.ActiveDocument

# 1. Capitalize "entwickeln" -> "Entwickeln" in the
#    "Ich interessiere mich sehr stark ..." paragraph.
$d.Content.Find.Execute("für das entwickeln von Web", $true, $false, $false, $false, $false, $true, 1, $false, "für das Entwickeln von Web", 2)

# 2. Relocate the _GoBack bookmark: it currently sits alone in the
#    empty paragraph right after "Ausbildung: ...". It needs to move
#    into the middle of the "Die Diplomarbeit basiert..." paragraph,
#    right after "weiterentwickeln " (splitting that run in two).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$r = $d.Content
$r.Find.Execute("weiterentwickeln ")
$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
